$rows = @(
    ,@(0,"falling",-0.9284301400184631,1.425136804580689,-0.2127189040184021,-0.09926560521125791,-0.0462730415165424,-0.0271835029125213)
    ,@(100,"falling",-0.9707106351852418,1.45836865901947,-0.1764526814222335,-0.0335975885391235,-0.0432187169790267,0.0745255574584007)
    ,@(200,"falling",-1.115207254886627,1.46594226360321,-0.1345747746527196,-0.0479529201984405,0.0383317954838275,0.0061086523346602)
    ,@(300,"falling",-1.505423545837401,1.455123424530029,-0.2340321838855745,0.046578474342823,-0.020616702735424,0.0142026171088218)
    ,@(400,"falling",-1.116380929946899,1.48697829246521,-0.4328413642942907,-0.0091629782691597,-0.0673478916287422,0.0209221355617046)
    ,@(500,"falling",-1.109515905380249,1.432106614112854,-0.3912773653864859,0.0207694191485643,-0.0343611687421798,0.0255036242306232)
    ,@(600,"falling",-1.141456544399262,1.384602665901184,-0.2541450988501308,0.022754730656743,-0.00534507073462,0.0320704244077205)
    ,@(700,"falling",-1.13429856300354,1.39785385131836,-0.2251825407147409,0.040775254368782,0.0120645882561802,0.009010262787342)
    ,@(800,"falling",-1.063723325729371,1.418689608573914,-0.2562501281499862,0.0372627787292003,-0.0259617734700441,0.0166460778564214)
    ,@(900,"falling",-1.018438935279846,1.406062006950379,-0.2269966453313826,0.011148290708661,-0.0271835029125213,0.0561996027827262)
    ,@(1000,"falling",-1.137969434261322,1.409385621547699,-0.1802991181612014,-0.0181732401251792,-0.0284052342176437,-0.011148290708661)
    ,@(1100,"falling",-1.03581714630127,1.381664276123047,-0.1787742376327515,0.0007635815418325,0.01328631862998,0.0305432621389627)
    ,@(1200,"falling",-1.100839495658875,1.406517148017883,-0.2175595723092557,0.0154243474826216,-0.0059559359215199,0.0029016099870204)
    ,@(1300,"falling",-1.232075214385986,1.368059515953064,-0.2088889628648757,0.0140499006956815,0.0103847095742821,0.0612392425537109)
    ,@(1400,"falling",-1.198645412921905,1.359035015106201,-0.2234921492636203,0.0195476878434419,0.0311541277915239,0.0459676086902618)
    ,@(1500,"falling",-1.162086248397827,1.317891120910645,0.1167446374893195,0.0041233403608202,-0.0007635815418325,-0.066737025976181)
    ,@(1600,"falling",-1.455766379833223,1.126043200492858,0.7863338142633457,-0.0314595587551593,0.2884811162948608,0.06536258012056349)
    ,@(1700,"falling",-1.653480172157284,0.8844107389450073,1.495913922786714,-0.1331686228513717,0.493731826543808,-0.1162171140313148)
    ,@(1800,"falling",-1.116286456584938,0.6627160906791623,2.260028153657919,-0.09651670604944219,0.7289149761199951,-0.0355829000473022)
    ,@(1900,"falling",-3.20666265487671,-0.9849638938903851,3.828832626342773,-0.2052507251501083,0.5688682794570923,0.1585195362567901)
    ,@(2000,"falling",-3.492754817008973,-1.960709273815156,3.70500636100769,-0.18539759516716,-0.1214094683527946,-0.6217080950737)
    ,@(2100,"falling",-2.38550305366516,0.4144415855407741,-0.4261573851108604,0.8633053302764893,-0.9155342578887939,0.0097738439217209)
    ,@(2200,"falling",-2.568133831024171,3.730020523071294,-0.8093817904591537,-2.978273391723633,1.564731359481812,-3.266449213027954)
    ,@(2300,"falling",-1.976609468460079,6.545797109603887,2.459241539239894,1.346346974372864,1.957823157310486,1.074512004852295)
    ,@(2400,"falling",-1.197644114494325,5.058232277631741,2.450196892023069,0.2492330223321914,0.9819658994674684,-0.9859365224838256)
    ,@(2500,"falling",-1.255056142807007,0.5559926331043243,-1.602012172341346,0.0514653958380222,-0.3197879493236542,0.1852448880672454)
    ,@(2600,"falling",-0.5990372896194439,1.264416024088864,-0.3774302378296863,0.1069014146924018,-0.1577559560537338,0.1151480972766876)
    ,@(2700,"falling",-0.1129188537597673,2.17881894111634,-0.6807380914687996,0.7513642311096191,-1.080926060676575,-0.1253800988197326)
    ,@(2800,"falling",-0.4983874559402485,3.747065991163262,1.904177859425558,-0.155312493443489,0.2770273983478546,0.052381694316864)
    ,@(2900,"falling",-0.2239453792572,0.9145344123244095,-1.341536760330224,-0.0971275717020034,-0.6624833345413208,0.3292563557624817)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = $rows.Count
$numCols = 8

$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $r = $rows[$i]
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $r[$j]
    }
}

$ws.Range("A2:H31").Value = $arr

Write-Output "Wrote $numRows rows into A2:H31"
